# Auto-generated Excel COM-interop script
# Commit: "Update latest output (run 155)"
# Updates the optimisation_result workbook: new pump-schedule run results on the
# "Schedule" sheet (rows 2-5, dimension grows to A1:F5) and refreshed forecast/
# historical price + Pump_Status values on the "Detailed" sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule": replace rows 2-4 with new optimisation values and add new row 5 ---
$wsSchedule = $wb.Worksheets.Item("Schedule")

# Extend formatting (incl. the date/time number format on columns A & B) from row 4
# down into the newly added row 5 before writing values into it.
$wsSchedule.Range("A4:F4").Copy()
$wsSchedule.Range("A5:F5").PasteSpecial(-4122)

$scheduleValues = @(
    @{Cell="A2"; Value=46043}
    @{Cell="B2"; Value=46043.16666666666}
    @{Cell="C2"; Value=4}
    @{Cell="D2"; Value=15.12}
    @{Cell="E2"; Value=506.74358475}
    @{Cell="F2"; Value=33.51478735119048}
    @{Cell="A3"; Value=46043.29166666666}
    @{Cell="B3"; Value=46043.66666666666}
    @{Cell="C3"; Value=9}
    @{Cell="D3"; Value=34.02}
    @{Cell="E3"; Value=-221.8333455}
    @{Cell="F3"; Value=-6.520674470899471}
    @{Cell="A4"; Value=46043.875}
    @{Cell="B4"; Value=46044.08333333334}
    @{Cell="C4"; Value=5}
    @{Cell="D4"; Value=18.9}
    @{Cell="E4"; Value=641.3837625}
    @{Cell="F4"; Value=33.93564880952381}
    @{Cell="A5"; Value=46044.25}
    @{Cell="B5"; Value=46044.66666666666}
    @{Cell="C5"; Value=10}
    @{Cell="D5"; Value=37.8}
    @{Cell="E5"; Value=52.11522225000004}
    @{Cell="F5"; Value=1.378709583333334}
)

foreach ($chg in $scheduleValues) {
    $wsSchedule.Range($chg.Cell).Value = $chg.Value
}

# --- Sheet "Detailed": update forecast/historical Price values and Pump_Status flags ---
$wsDetailed = $wb.Worksheets.Item("Detailed")

$detailedValues = @(
    @{Cell="E10"; Value="OFF"}
    @{Cell="E11"; Value="OFF"}
    @{Cell="E12"; Value="OFF"}
    @{Cell="E13"; Value="OFF"}
    @{Cell="E14"; Value="OFF"}
    @{Cell="E15"; Value="OFF"}
    @{Cell="E44"; Value="ON"}
    @{Cell="B45"; Value=59.66383}
    @{Cell="E45"; Value="ON"}
    @{Cell="B46"; Value=57.09}
    @{Cell="E46"; Value="ON"}
    @{Cell="B47"; Value=61.91795}
    @{Cell="C47"; Value="historical"}
    @{Cell="E47"; Value="ON"}
    @{Cell="C48"; Value="historical"}
    @{Cell="E48"; Value="ON"}
    @{Cell="B49"; Value=76.63499}
    @{Cell="C49"; Value="historical"}
    @{Cell="E49"; Value="ON"}
    @{Cell="B50"; Value=73.2}
    @{Cell="E50"; Value="ON"}
    @{Cell="B51"; Value=66.36179}
    @{Cell="E51"; Value="ON"}
    @{Cell="B52"; Value=66.16028}
    @{Cell="E52"; Value="ON"}
    @{Cell="B53"; Value=66.29066}
    @{Cell="E53"; Value="ON"}
    @{Cell="B54"; Value=66.0025}
    @{Cell="B55"; Value=66.91710999999999}
    @{Cell="B56"; Value=73.2}
    @{Cell="B57"; Value=73.2}
    @{Cell="B59"; Value=73.2}
    @{Cell="B60"; Value=67.38898}
    @{Cell="B62"; Value=64.89}
    @{Cell="E62"; Value="ON"}
    @{Cell="B63"; Value=57.06}
    @{Cell="E63"; Value="ON"}
    @{Cell="B64"; Value=35.88}
    @{Cell="E64"; Value="ON"}
    @{Cell="E65"; Value="ON"}
    @{Cell="B66"; Value=-5.95032}
    @{Cell="B67"; Value=-6.38688}
    @{Cell="B68"; Value=-6.90384}
    @{Cell="B69"; Value=-7.84163}
    @{Cell="B70"; Value=-7.67981}
    @{Cell="B71"; Value=-9.029170000000001}
    @{Cell="B72"; Value=-10}
    @{Cell="B73"; Value=-5.58973}
    @{Cell="B74"; Value=-6.72804}
    @{Cell="B75"; Value=-7.42351}
    @{Cell="B76"; Value=-7.85989}
    @{Cell="B77"; Value=-5.95281}
    @{Cell="B79"; Value=-2.54301}
    @{Cell="B80"; Value=0.51}
    @{Cell="B81"; Value=-10}
    @{Cell="B82"; Value=-9.710129999999999}
    @{Cell="B83"; Value=-10}
    @{Cell="B84"; Value=-11.01}
    @{Cell="B85"; Value=-8.114129999999999}
    @{Cell="B86"; Value=-1.63893}
    @{Cell="B87"; Value=7.77743}
    @{Cell="B88"; Value=18.5855}
    @{Cell="B89"; Value=51.4753}
    @{Cell="B90"; Value=59.36649}
    @{Cell="E90"; Value="OFF"}
    @{Cell="B91"; Value=61.94424}
    @{Cell="E91"; Value="OFF"}
    @{Cell="B92"; Value=57.04367}
    @{Cell="E92"; Value="OFF"}
    @{Cell="B93"; Value=57.03541}
    @{Cell="E93"; Value="OFF"}
    @{Cell="B94"; Value=42.99245}
    @{Cell="E94"; Value="OFF"}
    @{Cell="E95"; Value="OFF"}
    @{Cell="B96"; Value=56.98}
    @{Cell="E96"; Value="OFF"}
    @{Cell="B97"; Value=48.92621}
    @{Cell="E97"; Value="OFF"}
)

foreach ($chg in $detailedValues) {
    $wsDetailed.Range($chg.Cell).Value = $chg.Value
}

